$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Matches_SOG: append 3 new match rows (430-432)
# ---------------------------------------------------------------------------
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

function Set-MatchRow {
    param($ws, $row, $uid, $dateUtc, $home, $away, $sogHome, $sogAway, $source)

    # uid looks numeric ("897728") but must stay stored as text, matching the
    # rest of column A. A leading apostrophe forces text entry; re-applying
    # the Normal style afterwards drops the quote-prefix number format that
    # Excel otherwise tags the cell with, so no stray style index is left.
    $ws.Range("A$row").Value = "'" + $uid
    $ws.Range("A$row").Style = "Normal"

    $ws.Range("B$row").Value = $dateUtc
    $ws.Range("C$row").Value = $home
    $ws.Range("D$row").Value = $away
    $ws.Range("E$row").Value = $sogHome
    $ws.Range("F$row").Value = $sogAway
    $ws.Range("G$row").Value = $source
}

Set-MatchRow $wsMatches 430 "897728" "2025-11-04T13:30:00" "Сибирь"   "Торпедо"  25 39 "khl_text"
Set-MatchRow $wsMatches 431 "897726" "2025-11-04T17:00:00" "СКА"      "Динамо М" 23 36 "khl_text"
Set-MatchRow $wsMatches 432 "897729" "2025-11-04T17:00:00" "ХК Сочи"  "Спартак"  29 33 "khl_text"

# ---------------------------------------------------------------------------
# 2) Shots_HA: refresh as_of_utc stamp for every team row, then update the
#    home/away on-goal totals for the six teams that played on 2025-11-04.
# ---------------------------------------------------------------------------
$wsHA = $wb.Worksheets.Item("Shots_HA")

$newStamp = "2025-11-04T17:00:00Z"

for ($r = 2; $r -le 23; $r++) {
    $wsHA.Range("D$r").Value = $newStamp
}

$wsHA.Range("F8").Value  = 21
$wsHA.Range("K8").Value  = 596
$wsHA.Range("L8").Value  = 673
$wsHA.Range("M8").Value  = 28.4
$wsHA.Range("N8").Value  = 32

$wsHA.Range("E15").Value = 24
$wsHA.Range("G15").Value = 791
$wsHA.Range("H15").Value = 809
$wsHA.Range("I15").Value = 33
$wsHA.Range("J15").Value = 33.7

$wsHA.Range("E18").Value = 19
$wsHA.Range("G18").Value = 511
$wsHA.Range("H18").Value = 733
$wsHA.Range("I18").Value = 26.9

$wsHA.Range("F19").Value = 14
$wsHA.Range("K19").Value = 489
$wsHA.Range("L19").Value = 507
$wsHA.Range("M19").Value = 34.9
$wsHA.Range("N19").Value = 36.2

$wsHA.Range("F20").Value = 26
$wsHA.Range("K20").Value = 911
$wsHA.Range("L20").Value = 837
$wsHA.Range("M20").Value = 35
$wsHA.Range("N20").Value = 32.2

$wsHA.Range("E22").Value = 18
$wsHA.Range("G22").Value = 541
$wsHA.Range("H22").Value = 585

# ---------------------------------------------------------------------------
# 3) Shots_Summary: same as_of_utc refresh plus the matching total/pg updates
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Shots_Summary")

for ($r = 2; $r -le 23; $r++) {
    $wsSummary.Range("D$r").Value = $newStamp
}

$wsSummary.Range("E8").Value  = 36
$wsSummary.Range("F8").Value  = 1078
$wsSummary.Range("G8").Value  = 1088
$wsSummary.Range("H8").Value  = 29.9
$wsSummary.Range("I8").Value  = 30.2

$wsSummary.Range("E15").Value = 39
$wsSummary.Range("F15").Value = 1259
$wsSummary.Range("G15").Value = 1296
$wsSummary.Range("H15").Value = 32.3

$wsSummary.Range("E18").Value = 40
$wsSummary.Range("F18").Value = 1099
$wsSummary.Range("G18").Value = 1383
$wsSummary.Range("I18").Value = 34.6

$wsSummary.Range("E19").Value = 38
$wsSummary.Range("F19").Value = 1343
$wsSummary.Range("G19").Value = 1171
$wsSummary.Range("H19").Value = 35.3
$wsSummary.Range("I19").Value = 30.8

$wsSummary.Range("E20").Value = 46
$wsSummary.Range("F20").Value = 1554
$wsSummary.Range("G20").Value = 1433
$wsSummary.Range("H20").Value = 33.8
$wsSummary.Range("I20").Value = 31.2

$wsSummary.Range("E22").Value = 36
$wsSummary.Range("F22").Value = 993
$wsSummary.Range("G22").Value = 1258
$wsSummary.Range("H22").Value = 27.6
$wsSummary.Range("I22").Value = 34.9

# ---------------------------------------------------------------------------
# 4) Meta_ext: bump the as_of timestamp and build_version
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Meta_ext")
$wsMeta.Range("B2").Value = $newStamp
$wsMeta.Range("D2").Value = 43
